# Daily attendance processing - 2025-11-03 08:29:45
#
# Column G ("Recorded By") holds a comma-separated list of the
# actors who recorded/touched each attendance session (real user
# emails, plus synthetic/system actors "System", "system" and
# "admin@admin.com"). The recorder list was being built with the
# synthetic actors sorted ahead of real users; this pass re-sorts
# each cell's list so real user identities sort first and the
# synthetic/system actors are pushed toward the end, in the order:
#   (real users) -> admin@admin.com -> System -> system
# Relative order within each priority group is preserved (stable
# sort), and lists with a single entry are left untouched.

function Get-ActorRank($item) {
    if ($item -eq "System") { return 2 }
    if ($item -eq "system") { return 3 }
    if ($item -eq "admin@admin.com") { return 1 }
    return 0
}

function Reorder-RecordedBy($val) {
    $parts = @($val -split ", ")
    if ($parts.Length -le 1) {
        return $val
    }

    $ranks = @()
    foreach ($p in $parts) {
        $ranks += (Get-ActorRank $p)
    }

    # Manual stable insertion sort (avoids pipeline quirks of
    # Sort-Object when operating on tiny / single-element arrays).
    $n = $parts.Length
    for ($i = 1; $i -lt $n; $i++) {
        $keyPart = $parts[$i]
        $keyRank = $ranks[$i]
        $j = $i - 1
        while ($j -ge 0 -and $ranks[$j] -gt $keyRank) {
            $parts[$j + 1] = $parts[$j]
            $ranks[$j + 1] = $ranks[$j]
            $j = $j - 1
        }
        $parts[$j + 1] = $keyPart
        $ranks[$j + 1] = $keyRank
    }

    return ($parts -join ", ")
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $current = $cell.Value2
    if ($current -ne $null -and $current -ne "") {
        $updated = Reorder-RecordedBy $current
        if ($updated -ne $current) {
            $cell.Value = $updated
        }
    }
}
